$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 240
$ws.Range("J55").Value = 287.33334
$ws.Range("L55").Value = 287.33334
$ws.Range("N55").Value = -715.33334
$ws.Range("H100").Value = 2405.8125
$ws.Range("I100").Value = 591.5
$ws.Range("K100").Value = 591.5
$ws.Range("M100").Value = -50.5
$ws.Range("H111").Value = 562.2143
$ws.Range("I111").Value = 547.5833
$ws.Range("K111").Value = 1642.7499
$ws.Range("M111").Value = 1424.2501
$ws.Range("H116").Value = 5664.3335
$ws.Range("I116").Value = 5664.3335
$ws.Range("K116").Value = 5664.3335
$ws.Range("M116").Value = -2222.3335
$ws.Range("H135").Value = 792.5714
$ws.Range("I135").Value = 752
$ws.Range("K135").Value = 6768
$ws.Range("M135").Value = -4233

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2428
$ws.Range("I74").Value = 3856
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 3856
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -2982
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 2428
$ws.Range("I77").Value = 3856
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 19280
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -14912
$ws.Range("N77").Value = -13736

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I94").Value = 2622
$ws.Range("J94").Value = 1960
$ws.Range("K94").Value = 2622
$ws.Range("L94").Value = 1960
$ws.Range("M94").Value = -2171
$ws.Range("N94").Value = -2862
$ws.Range("H134").Value = 5410.5
$ws.Range("I134").Value = 5410.5
$ws.Range("K134").Value = 16231.5
$ws.Range("M134").Value = -13696.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2396.5
$ws.Range("I31").Value = 2396.5
$ws.Range("K31").Value = 2396.5
$ws.Range("M31").Value = -2101.5
$ws.Range("H34").Value = 2396.5
$ws.Range("I34").Value = 2396.5
$ws.Range("K34").Value = 2396.5
$ws.Range("M34").Value = -2194.5
$ws.Range("H43").Value = 25552.334
$ws.Range("J43").Value = 25552.334
$ws.Range("L43").Value = 25552.334
$ws.Range("N43").Value = -25920.334
$ws.Range("H99").Value = 4577.8
$ws.Range("I99").Value = 3949.5
$ws.Range("J99").Value = 4996.6665
$ws.Range("K99").Value = 3949.5
$ws.Range("L99").Value = 4996.6665
$ws.Range("M99").Value = -2451.5
$ws.Range("N99").Value = -7992.6665
$ws.Range("H101").Value = 25552.334
$ws.Range("J101").Value = 25552.334
$ws.Range("L101").Value = 25552.334
$ws.Range("N101").Value = -32042.334
$ws.Range("H126").Value = 4577.8
$ws.Range("I126").Value = 3949.5
$ws.Range("J126").Value = 4996.6665
$ws.Range("K126").Value = 11848.5
$ws.Range("L126").Value = 14989.9995
$ws.Range("M126").Value = -9378.5
$ws.Range("N126").Value = -19929.9995
$ws.Range("H134").Value = 2501
$ws.Range("I134").Value = 2626.0625
$ws.Range("J134").Value = 500
$ws.Range("K134").Value = 7878.1875
$ws.Range("L134").Value = 1500
$ws.Range("M134").Value = -5343.1875
$ws.Range("N134").Value = -6570

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1538.4546
$ws.Range("J5").Value = 489.5
$ws.Range("L5").Value = 1468.5
$ws.Range("N5").Value = -1692.5
$ws.Range("H12").Value = 179.36363
$ws.Range("I12").Value = 42.75
$ws.Range("J12").Value = 257.42856
$ws.Range("K12").Value = 128.25
$ws.Range("L12").Value = 772.28568
$ws.Range("M12").Value = 44.75
$ws.Range("N12").Value = -1118.28568
$ws.Range("H61").Value = 460.375
$ws.Range("I61").Value = 363.83334
$ws.Range("J61").Value = 750
$ws.Range("K61").Value = 1091.50002
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -876.5000199999999
$ws.Range("N61").Value = -2680
$ws.Range("H68").Value = 2002.5
$ws.Range("J68").Value = 2002.5
$ws.Range("L68").Value = 6007.5
$ws.Range("N68").Value = -7629.5
$ws.Range("H71").Value = 2002.5
$ws.Range("J71").Value = 2002.5
$ws.Range("L71").Value = 18022.5
$ws.Range("N71").Value = -26134.5
$ws.Range("H131").Value = 1540.3334
$ws.Range("J131").Value = 1998
$ws.Range("L131").Value = 5994
$ws.Range("N131").Value = -16074
$ws.Range("H135").Value = 1538.4546
$ws.Range("J135").Value = 489.5
$ws.Range("L135").Value = 4405.5
$ws.Range("N135").Value = -9475.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 5000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984
$ws.Range("H101").Value = 35999.5
$ws.Range("J101").Value = 35999.5
$ws.Range("L101").Value = 35999.5
$ws.Range("N101").Value = -42489.5
$ws.Range("H113").Value = 3292.2
$ws.Range("I113").Value = 3292.2
$ws.Range("K113").Value = 3292.2
$ws.Range("M113").Value = -1122.2
$ws.Range("H122").Value = 11367206
$ws.Range("I122").Value = 17859880
$ws.Range("K122").Value = 53579640
$ws.Range("M122").Value = -53577190

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10272
$ws.Range("H68").Value = 2049.1428
$ws.Range("I68").Value = 2119
$ws.Range("J68").Value = 1874.5
$ws.Range("K68").Value = 2119
$ws.Range("L68").Value = 1874.5
$ws.Range("M68").Value = -1370
$ws.Range("N68").Value = -3372.5
$ws.Range("H71").Value = 2049.1428
$ws.Range("I71").Value = 2119
$ws.Range("J71").Value = 1874.5
$ws.Range("K71").Value = 10595
$ws.Range("L71").Value = 9372.5
$ws.Range("M71").Value = -6851
$ws.Range("N71").Value = -16860.5
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H101").Value = 18135.625
$ws.Range("J101").Value = 18135.625
$ws.Range("L101").Value = 18135.625
$ws.Range("N101").Value = -24625.625
$ws.Range("M40").ClearContents()
$ws.Range("N93").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4490.2666
$ws.Range("I62").Value = 4696.5835
$ws.Range("J62").Value = 3665
$ws.Range("K62").Value = 4696.5835
$ws.Range("L62").Value = 3665
$ws.Range("M62").Value = -4072.5835
$ws.Range("N62").Value = -4913
$ws.Range("H65").Value = 4490.2666
$ws.Range("I65").Value = 4696.5835
$ws.Range("J65").Value = 3665
$ws.Range("K65").Value = 23482.9175
$ws.Range("L65").Value = 18325
$ws.Range("M65").Value = -20362.9175
$ws.Range("N65").Value = -24565
$ws.Range("H107").Value = 579.1429000000001
$ws.Range("I107").Value = 579.1429000000001
$ws.Range("K107").Value = 1737.4287
$ws.Range("M107").Value = 182.5712999999998
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330
